$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $escaped = $val -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextValue "E2" "2026-02-06 05:48:01"
Set-TextValue "N2" "-2.7 °C 5:14 TU"
Set-TextValue "E3" "2026-02-06 05:48:03"
Set-TextValue "K3" "-0.1 MJ/m2"
Set-TextValue "E4" "2026-02-06 05:48:06"
Set-TextValue "H4" "60%"
Set-TextValue "J4" "993.3 hPa"
Set-TextValue "N4" "8.0 °C 5:28 TU"
Set-TextValue "O4" "12.4 °C"
Set-TextValue "E5" "2026-02-06 05:48:09"
Set-TextValue "J5" "993.8 hPa"
Set-TextValue "O5" "8.0 °C"
Set-TextValue "E6" "2026-02-06 05:48:11"
Set-TextValue "H6" "51%"
Set-TextValue "J6" "994.9 hPa"
Set-TextValue "N6" "13.5 °C 5:29 TU"
Set-TextValue "E7" "2026-02-06 05:48:13"
Set-TextValue "J7" "994.7 hPa"
Set-TextValue "N7" "9.0 °C 5:29 TU"
Set-TextValue "E8" "2026-02-06 05:48:16"
Set-TextValue "H8" "93%"
Set-TextValue "N8" "4.1 °C 5:00 TU"
Set-TextValue "O8" "5.8 °C"
Set-TextValue "E9" "2026-02-06 05:48:18"
Set-TextValue "N9" "0.5 °C 5:24 TU"
Set-TextValue "O9" "2.1 °C"
Set-TextValue "E10" "2026-02-06 05:48:20"
Set-TextValue "O10" "5.0 °C"
Set-TextValue "E11" "2026-02-06 05:48:23"
Set-TextValue "J11" "995.6 hPa"
Set-TextValue "N11" "1.7 °C 5:23 TU"
Set-TextValue "O11" "4.4 °C"
Set-TextValue "E12" "2026-02-06 05:48:26"
Set-TextValue "H12" "61%"
Set-TextValue "O12" "12.3 °C"
Set-TextValue "E13" "2026-02-06 05:48:28"
Set-TextValue "H13" "90%"
Set-TextValue "N13" "3.3 °C 5:29 TU"
Set-TextValue "O13" "6.2 °C"
Set-TextValue "E14" "2026-02-06 05:48:30"
Set-TextValue "H14" "73%"
Set-TextValue "E15" "2026-02-06 05:48:33"
Set-TextValue "H15" "86%"
Set-TextValue "J15" "993.9 hPa"
Set-TextValue "N15" "3.0 °C 5:29 TU"
Set-TextValue "O15" "6.6 °C"
Set-TextValue "E16" "2026-02-06 05:48:36"
Set-TextValue "N16" "3.1 °C 5:11 TU"
Set-TextValue "O16" "4.1 °C"
Set-TextValue "E17" "2026-02-06 05:48:39"
Set-TextValue "J17" "997.0 hPa"
Set-TextValue "N17" "1.0 °C 5:28 TU"
Set-TextValue "O17" "3.0 °C"
Set-TextValue "E18" "2026-02-06 05:48:41"
Set-TextValue "O18" "-4.9 °C"
Set-TextValue "E19" "2026-02-06 05:48:44"
Set-TextValue "H19" "98%"
Set-TextValue "J19" "997.3 hPa"
Set-TextValue "E20" "2026-02-06 05:48:47"
Set-TextValue "K20" "-0.1 MJ/m2"
Set-TextValue "O20" "-2.2 °C"
Set-TextValue "E21" "2026-02-06 05:48:50"
Set-TextValue "J21" "994.8 hPa"
Set-TextValue "N21" "2.7 °C 5:04 TU"
Set-TextValue "O21" "4.8 °C"
Set-TextValue "E22" "2026-02-06 05:48:52"
Set-TextValue "H22" "83%"
Set-TextValue "N22" "3.8 °C 5:12 TU"
Set-TextValue "O22" "7.7 °C"
Set-TextValue "E23" "2026-02-06 05:48:55"
Set-TextValue "J23" "994.0 hPa"
Set-TextValue "N23" "6.4 °C 5:04 TU"
Set-TextValue "E24" "2026-02-06 05:48:58"
Set-TextValue "J24" "992.8 hPa"
Set-TextValue "E25" "2026-02-06 05:49:01"
Set-TextValue "H25" "93%"
Set-TextValue "J25" "996.1 hPa"
Set-TextValue "O25" "2.0 °C"
Set-TextValue "E26" "2026-02-06 05:49:03"
Set-TextValue "H26" "80%"
Set-TextValue "N26" "-3.1 °C 5:29 TU"
Set-TextValue "O26" "-0.5 °C"
Set-TextValue "E27" "2026-02-06 05:49:05"
Set-TextValue "J27" "993.7 hPa"
Set-TextValue "O27" "7.3 °C"
Set-TextValue "E28" "2026-02-06 05:49:08"
Set-TextValue "H28" "90%"
Set-TextValue "J28" "996.8 hPa"
Set-TextValue "N28" "-0.3 °C 5:20 TU"
Set-TextValue "O28" "2.6 °C"
Set-TextValue "E29" "2026-02-06 05:49:10"
Set-TextValue "H29" "65%"
Set-TextValue "N29" "6.5 °C 5:29 TU"
Set-TextValue "O29" "11.1 °C"
Set-TextValue "E30" "2026-02-06 05:49:13"
Set-TextValue "H30" "75%"
Set-TextValue "K30" "-0.1 MJ/m2"
Set-TextValue "E31" "2026-02-06 05:49:16"
Set-TextValue "J31" "996.9 hPa"
Set-TextValue "O31" "4.9 °C"
Set-TextValue "E32" "2026-02-06 05:49:18"
Set-TextValue "J32" "995.3 hPa"
Set-TextValue "O32" "14.7 °C"
Set-TextValue "E33" "2026-02-06 05:49:21"
Set-TextValue "N33" "5.0 °C 5:29 TU"
Set-TextValue "O33" "6.6 °C"
Set-TextValue "E34" "2026-02-06 05:49:23"
Set-TextValue "H34" "80%"
Set-TextValue "K34" "-0.1 MJ/m2"
Set-TextValue "N34" "3.1 °C 5:00 TU"
Set-TextValue "O34" "7.5 °C"
Set-TextValue "E35" "2026-02-06 05:49:26"
Set-TextValue "N35" "-3.4 °C 5:13 TU"
Set-TextValue "E36" "2026-02-06 05:49:29"
Set-TextValue "J36" "996.7 hPa"
Set-TextValue "O36" "11.4 °C"
